# freq.xlsx: add "disturbance" (column A) labels, abbreviate "obs" (column C)
# values, and fix a species-code typo ("jute" -> "juge") in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: disturbance group labels -----------------------------------
# The data block (rows 2-94) is divided into four contiguous disturbance
# groups. Fill each block in one shot.
$ws.Range("A2:A27").Value  = "Grubbed"
$ws.Range("A28:A47").Value = "Heavily Grazed, `nExclosed 1 year"
$ws.Range("A48:A62").Value = "Heavily Grazed, `nExclosed 10 years"
$ws.Range("A63:A94").Value = "Undisturbed"

# --- Columns B & C: fix typo + abbreviate observation type -----------------
for ($r = 2; $r -le 94; $r++) {
    $c = $ws.Cells.Item($r, 3).Value2
    if ($c -eq "Surface seed bank") {
        $ws.Cells.Item($r, 3).Value = "seed"
    } elseif ($c -eq "Above-ground vegetation") {
        $ws.Cells.Item($r, 3).Value = "veg"
    }

    $b = $ws.Cells.Item($r, 2).Value2
    if ($b -eq "jute") {
        $ws.Cells.Item($r, 2).Value = "juge"
    }
}
